$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,('MS Dhoni (WK)', 1, 0, 0, 28, 50)
    ,('Mustafizur Rahman', 0, 0, 0, 24, 29)
    ,('Virat Kohli', 1, 0, 0, 18, 27)
    ,('Tushar Deshpande', 0, 0, 0, 23, 27)
    ,('Deepak Chahar', 1, 0, 0, 17, 27)
    ,('Cameron Green', 2, 0, 0, 16, 25)
    ,('Ruturaj Gaikwad (C)', 0, 0, 0, 14, 24)
    ,('Ravindra Jadeja', 1, 0, 0, 11, 24)
    ,('Alzarri Josheph', 2, 0, 0, 13, 23)
    ,('Faf du Plessis(C)', 1, 0, 0, 18, 22)
    ,('Mayank Dagar', 2, 0, 0, 6, 21)
    ,('Mohammed Siraj', 2, 0, 0, 10, 20)
    ,('Rachin Ravindra', 0, 0, 0, 7, 19)
    ,('Dinesh Karthik', 0, 0, 0, 10, 18)
    ,('Maheesh Theekshana', 0, 0, 0, 10, 17)
    ,('Karn Sharma', 1, 0, 0, 4, 17)
    ,('Anuj Rawat(WK)', 0, 0, 0, 13, 17)
    ,('Sameer Rizvi', 1, 0, 0, 12, 16)
    ,('Daryl Mitchell', 0, 0, 0, 2, 13)
    ,('Ajinkya Rahane', 0, 0, 0, 8, 10)
    ,('Glenn Maxwell', 0, 0, 0, 4, 7)
    ,('Rajat Patidar', 0, 0, 0, 4, 5)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $item = $data[$i]
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $ws.Cells.Item($row, 6).Value = $item[5]
}
